$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, week-of dates) ---
$ws.Cells.Item(8, 1).Value = "Volume 30   Number  42"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  10/16/2023  Through  10/22/2023"

# --- Plain numeric value updates (precinct weekly crime stats) ---
$ws.Cells.Item(14, 14).Value = -62.962962962963
$ws.Cells.Item(15, 6).Value = 2
$ws.Cells.Item(15, 8).Value = -50
$ws.Cells.Item(15, 9).Value = 25
$ws.Cells.Item(15, 11).Value = -26.470588235294
$ws.Cells.Item(15, 12).Value = -19.354838709677
$ws.Cells.Item(15, 13).Value = 92.307692307692
$ws.Cells.Item(15, 14).Value = -62.121212121212
$ws.Cells.Item(16, 3).Value = 9
$ws.Cells.Item(16, 5).Value = 50
$ws.Cells.Item(16, 6).Value = 38
$ws.Cells.Item(16, 7).Value = 25
$ws.Cells.Item(16, 8).Value = 52
$ws.Cells.Item(16, 9).Value = 345
$ws.Cells.Item(16, 10).Value = 346
$ws.Cells.Item(16, 11).Value = -0.289017341040
$ws.Cells.Item(16, 12).Value = 28.252788104089
$ws.Cells.Item(16, 13).Value = 38
$ws.Cells.Item(16, 14).Value = -66.176470588235
$ws.Cells.Item(17, 3).Value = 12
$ws.Cells.Item(17, 4).Value = 13
$ws.Cells.Item(17, 5).Value = -7.692307692307
$ws.Cells.Item(17, 6).Value = 38
$ws.Cells.Item(17, 7).Value = 43
$ws.Cells.Item(17, 8).Value = -11.627906976744
$ws.Cells.Item(17, 9).Value = 570
$ws.Cells.Item(17, 10).Value = 540
$ws.Cells.Item(17, 11).Value = 5.555555555555
$ws.Cells.Item(17, 12).Value = 8.778625954198
$ws.Cells.Item(17, 13).Value = 119.230769230769
$ws.Cells.Item(17, 14).Value = -28.571428571428
$ws.Cells.Item(18, 3).Value = 11
$ws.Cells.Item(18, 4).Value = 2
$ws.Cells.Item(18, 5).Value = 450
$ws.Cells.Item(18, 6).Value = 18
$ws.Cells.Item(18, 7).Value = 21
$ws.Cells.Item(18, 8).Value = -14.285714285714
$ws.Cells.Item(18, 9).Value = 163
$ws.Cells.Item(18, 10).Value = 272
$ws.Cells.Item(18, 11).Value = -40.073529411764
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = 29.365079365079
$ws.Cells.Item(18, 14).Value = -81.767337807606
$ws.Cells.Item(19, 4).Value = 7
$ws.Cells.Item(19, 5).Value = 42.857142857142
$ws.Cells.Item(19, 6).Value = 42
$ws.Cells.Item(19, 7).Value = 50
$ws.Cells.Item(19, 8).Value = -16
$ws.Cells.Item(19, 9).Value = 416
$ws.Cells.Item(19, 10).Value = 424
$ws.Cells.Item(19, 11).Value = -1.886792452830
$ws.Cells.Item(19, 12).Value = 8.051948051948
$ws.Cells.Item(19, 13).Value = 111.167512690355
$ws.Cells.Item(19, 14).Value = 39.130434782608
$ws.Cells.Item(20, 3).Value = 5
$ws.Cells.Item(20, 5).Value = -16.666666666666
$ws.Cells.Item(20, 6).Value = 27
$ws.Cells.Item(20, 7).Value = 25
$ws.Cells.Item(20, 8).Value = 8
$ws.Cells.Item(20, 9).Value = 362
$ws.Cells.Item(20, 10).Value = 248
$ws.Cells.Item(20, 11).Value = 45.967741935483
$ws.Cells.Item(20, 12).Value = 101.111111111111
$ws.Cells.Item(20, 13).Value = 289.247311827957
$ws.Cells.Item(20, 14).Value = -5.235602094240
$ws.Cells.Item(21, 3).Value = 48
$ws.Cells.Item(21, 4).Value = 34
$ws.Cells.Item(21, 5).Value = 41.176470588235
$ws.Cells.Item(21, 6).Value = 166
$ws.Cells.Item(21, 7).Value = 168
$ws.Cells.Item(21, 8).Value = -1.190476190476
$ws.Cells.Item(21, 9).Value = 1891
$ws.Cells.Item(21, 10).Value = 1869
$ws.Cells.Item(21, 11).Value = 1.177100053504
$ws.Cells.Item(21, 12).Value = 20.522625876354
$ws.Cells.Item(21, 13).Value = 99.472573839662
$ws.Cells.Item(21, 14).Value = -45.754446356856
$ws.Cells.Item(22, 8).Value = -50
$ws.Cells.Item(22, 9).Value = 9
$ws.Cells.Item(22, 11).Value = -18.181818181818
$ws.Cells.Item(22, 12).Value = -10
$ws.Cells.Item(22, 13).Value = -40
$ws.Cells.Item(23, 3).Value = 7
$ws.Cells.Item(23, 4).Value = 10
$ws.Cells.Item(23, 5).Value = -30
$ws.Cells.Item(23, 6).Value = 25
$ws.Cells.Item(23, 7).Value = 29
$ws.Cells.Item(23, 8).Value = -13.793103448275
$ws.Cells.Item(23, 9).Value = 332
$ws.Cells.Item(23, 10).Value = 292
$ws.Cells.Item(23, 11).Value = 13.698630136986
$ws.Cells.Item(23, 12).Value = 88.636363636363
$ws.Cells.Item(23, 13).Value = 115.584415584416
$ws.Cells.Item(24, 3).Value = 30
$ws.Cells.Item(24, 4).Value = 24
$ws.Cells.Item(24, 5).Value = 25
$ws.Cells.Item(24, 6).Value = 91
$ws.Cells.Item(24, 7).Value = 81
$ws.Cells.Item(24, 8).Value = 12.345679012345
$ws.Cells.Item(24, 9).Value = 966
$ws.Cells.Item(24, 10).Value = 1031
$ws.Cells.Item(24, 11).Value = -6.304558680892
$ws.Cells.Item(24, 12).Value = 25.454545454545
$ws.Cells.Item(24, 13).Value = 52.60663507109
$ws.Cells.Item(25, 3).Value = 18
$ws.Cells.Item(25, 4).Value = 9
$ws.Cells.Item(25, 5).Value = 100
$ws.Cells.Item(25, 6).Value = 62
$ws.Cells.Item(25, 7).Value = 64
$ws.Cells.Item(25, 8).Value = -3.125
$ws.Cells.Item(25, 9).Value = 868
$ws.Cells.Item(25, 10).Value = 802
$ws.Cells.Item(25, 11).Value = 8.229426433915
$ws.Cells.Item(25, 12).Value = 19.067215363511
$ws.Cells.Item(25, 13).Value = 21.229050279329
$ws.Cells.Item(26, 3).Value = 2
$ws.Cells.Item(26, 9).Value = 45
$ws.Cells.Item(26, 11).Value = -15.094339622641
$ws.Cells.Item(26, 12).Value = -8.163265306122
$ws.Cells.Item(27, 3).Value = 4
$ws.Cells.Item(27, 4).Value = 2
$ws.Cells.Item(27, 5).Value = 100
$ws.Cells.Item(27, 6).Value = 10
$ws.Cells.Item(27, 8).Value = 11.111111111111
$ws.Cells.Item(27, 9).Value = 84
$ws.Cells.Item(27, 10).Value = 65
$ws.Cells.Item(27, 11).Value = 29.230769230769
$ws.Cells.Item(27, 12).Value = 55.555555555555
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 12).Value = -47.692307692307
$ws.Cells.Item(28, 14).Value = -63.043478260869
$ws.Cells.Item(29, 6).Value = 1
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 12).Value = -49.090909090909
$ws.Cells.Item(29, 14).Value = -69.230769230769

# --- Cells that change between numeric and placeholder-text ("0" / "***.*") ---
# Use PasteSpecial(xlPasteFormats) from a same-style template cell so the target
# cell picks up the correct number format / style index before the value is written,
# matching how Excel itself keeps the "no data" placeholder cells styled.
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4122)
$ws.Cells.Item(15, 4).Value = "0"
$ws.Cells.Item(14, 5).Copy()
$ws.Cells.Item(15, 5).PasteSpecial(-4122)
$ws.Cells.Item(15, 5).Value = "***.*"
$ws.Cells.Item(14, 6).Copy()
$ws.Cells.Item(22, 3).PasteSpecial(-4122)
$ws.Cells.Item(22, 3).Value = 1
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4122)
$ws.Cells.Item(22, 4).Value = "0"
$ws.Cells.Item(14, 5).Copy()
$ws.Cells.Item(22, 5).PasteSpecial(-4122)
$ws.Cells.Item(22, 5).Value = "***.*"
$ws.Cells.Item(14, 6).Copy()
$ws.Cells.Item(22, 6).PasteSpecial(-4122)
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item(26, 4).PasteSpecial(-4122)
$ws.Cells.Item(26, 4).Value = "0"
$ws.Cells.Item(14, 5).Copy()
$ws.Cells.Item(26, 5).PasteSpecial(-4122)
$ws.Cells.Item(26, 5).Value = "***.*"

$excel.CutCopyMode = 0
